$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds numeric-looking text (e.g. "508.90", "1.00").
# Force those specific cells to Text format first so assigning the string
# value keeps it verbatim instead of Excel coercing it to a Number (which
# would silently drop meaningful trailing/grouping zeros, e.g. "508.90" -> 508.9).
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"

$ws.Range('D2').Value = '60.433.10'
$ws.Range('E2').Value = '  -1.54%  '
$ws.Range('D3').Value = '2.591.68'
$ws.Range('E3').Value = '  -2.94%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '508.90'
$ws.Range('E5').Value = '  -1.31%  '
$ws.Range('D6').Value = '154.66'
$ws.Range('E6').Value = '  -3.53%  '
$ws.Range('E7').Value = '  -0.15%  '
$ws.Range('D8').Value = '0.582'
$ws.Range('E8').Value = '  -5.17%  '
$ws.Range('D9').Value = '2.601.59'
$ws.Range('E9').Value = '  -2.60%  '
$ws.Range('E10').Value = '  +8.52%  '
$ws.Range('E11').Value = '  -2.07%  '
$ws.Range('E12').Value = '  -0.94%  '
$ws.Range('E13').Value = '  +1.19%  '
$ws.Range('D14').Value = '3.049.36'
$ws.Range('E14').Value = '  -2.55%  '
$ws.Range('D15').Value = '60.463.87'
$ws.Range('E15').Value = '  -1.45%  '
$ws.Range('D16').Value = '21.63'
$ws.Range('E16').Value = '  -3.58%  '
$ws.Range('E17').Value = '  -0.66%  '
$ws.Range('D18').Value = '2.601.01'
$ws.Range('E18').Value = '  -2.56%  '
$ws.Range('D19').Value = '4.76'
$ws.Range('E19').Value = '  -1.80%  '
$ws.Range('D20').Value = '348.19'
$ws.Range('E20').Value = '  -1.37%  '
$ws.Range('D21').Value = '10.53'
$ws.Range('E21').Value = '  -0.49%  '
$ws.Range('D22').Value = '6.13'
$ws.Range('E22').Value = '  -1.49%  '
$ws.Range('D23').Value = '0.997'
$ws.Range('E23').Value = '  -0.16%  '
$ws.Range('D24').Value = '60.55'
$ws.Range('E24').Value = '  -0.40%  '
$ws.Range('E25').Value = '  -1.63%  '
$ws.Range('E26').Value = '  -0.80%  '
$ws.Range('B27').Value = 'Binance-PegBSC-USD'
$ws.Range('C27').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D27').Value = '0.999'
$ws.Range('E27').Value = '  +0.17%  '
$ws.Range('B28').Value = 'PEPE'
$ws.Range('C28').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D28').Value = '0.0₃0846'
$ws.Range('E28').Value = '  -3.77%  '
$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').Value = '7.42'
$ws.Range('E29').Value = '  -1.74%  '
$ws.Range('B30').Value = 'USDe'
$ws.Range('C30').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').Value = '  -0.04%  '
$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').Value = '19.40'
$ws.Range('E31').Value = '  -1.87%  '
$ws.Range('B32').Value = 'Monero'
$ws.Range('C32').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D32').Value = '152.53'
$ws.Range('E32').Value = '  -3.03%  '
$ws.Range('B33').Value = 'PancakeSwap'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D33').Value = '1.56'
$ws.Range('E33').Value = '  -2.05%  '
$ws.Range('B34').Value = 'Aptos'
$ws.Range('C34').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D34').Value = '5.74'
$ws.Range('E34').Value = '  -0.70%  '
$ws.Range('B35').Value = 'NEARProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D35').Value = '4.04'
$ws.Range('E35').Value = '  -2.02%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').Value = '1.19'
$ws.Range('E36').Value = '  -3.74%  '
$ws.Range('B37').Value = 'SuiNetwork'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D37').Value = '0.857'
$ws.Range('E37').Value = '  +2.04%  '
$ws.Range('B38').Value = 'Stacks'
$ws.Range('C38').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D38').Value = '1.49'
$ws.Range('E38').Value = '  -4.74%  '
$ws.Range('B39').Value = 'Fetch.AI'
$ws.Range('C39').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D39').Value = '0.849'
$ws.Range('E39').Value = '  -3.67%  '
$ws.Range('B40').Value = 'Filecoin'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D40').Value = '3.78'
$ws.Range('E40').Value = '  -0.50%  '
$ws.Range('B41').Value = 'OKB'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D41').Value = '36.20'
$ws.Range('E41').Value = '  +1.11%  '
$ws.Range('B42').Value = 'Bittensor'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D42').Value = '298.75'
$ws.Range('E42').Value = '  -1.57%  '
$ws.Range('B43').Value = 'Mantle'
$ws.Range('C43').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D43').Value = '0.623'
$ws.Range('E43').Value = '  -3.89%  '
$ws.Range('B44').Value = 'Stellar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D44').Value = '0.1000'
$ws.Range('E44').Value = '  -2.28%  '
$ws.Range('B45').Value = 'Hedera'
$ws.Range('C45').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D45').Value = '0.0558'
$ws.Range('E45').Value = '  -3.80%  '
$ws.Range('B46').Value = 'FirstDigitalUSD'
$ws.Range('C46').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D46').Value = '0.998'
$ws.Range('E46').Value = '  -0.18%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = '19.85'
$ws.Range('E47').Value = '  -2.01%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').Value = '4.83'
$ws.Range('E48').Value = '  -3.51%  '
$ws.Range('B49').Value = 'VeChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D49').Value = '0.0234'
$ws.Range('E49').Value = '  -2.15%  '
$ws.Range('B50').Value = 'WhiteBITCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D50').Value = '10.30'
$ws.Range('E50').Value = '  -0.44%  '
$ws.Range('B51').Value = 'Maker'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D51').Value = '1.994.95'
$ws.Range('E51').Value = '  -2.02%  '
